$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1715
$ws.Range("J17").Value = 1835.625
$ws.Range("L17").Value = 5506.875
$ws.Range("N17").Value = -5842.875
$ws.Range("H62").Value = 5009.125
$ws.Range("I62").Value = 5029.7334
$ws.Range("K62").Value = 5029.7334
$ws.Range("M62").Value = -4405.7334
$ws.Range("H65").Value = 5009.125
$ws.Range("I65").Value = 5029.7334
$ws.Range("K65").Value = 25148.667
$ws.Range("M65").Value = -22028.667
$ws.Range("H80").Value = 3365.5518
$ws.Range("I80").Value = 1217.9
$ws.Range("J80").Value = 4495.8945
$ws.Range("K80").Value = 3653.7
$ws.Range("L80").Value = 13487.6835
$ws.Range("M80").Value = -2655.7
$ws.Range("N80").Value = -15483.6835
$ws.Range("H83").Value = 3365.5518
$ws.Range("I83").Value = 1217.9
$ws.Range("J83").Value = 4495.8945
$ws.Range("K83").Value = 10961.1
$ws.Range("L83").Value = 40463.0505
$ws.Range("M83").Value = -5969.1
$ws.Range("N83").Value = -50447.0505
$ws.Range("H92").Value = 3788918.5
$ws.Range("J92").Value = 9616139
$ws.Range("L92").Value = 9616139
$ws.Range("N92").Value = -9618635
$ws.Range("H98").Value = 1112090.8
$ws.Range("I98").Value = 1235100.9
$ws.Range("K98").Value = 1235100.9
$ws.Range("M98").Value = -1233602.9
$ws.Range("H101").Value = 989.9375
$ws.Range("I101").Value = 403.66666
$ws.Range("J101").Value = 2748.75
$ws.Range("K101").Value = 1210.99998
$ws.Range("L101").Value = 8246.25
$ws.Range("M101").Value = 411.0000199999999
$ws.Range("N101").Value = -11490.25
$ws.Range("H106").Value = 3549
$ws.Range("I106").Value = 3549
$ws.Range("K106").Value = 3549
$ws.Range("M106").Value = -2918
$ws.Range("H113").Value = 12000
$ws.Range("I113").Value = 15000
$ws.Range("K113").Value = 15000
$ws.Range("M113").Value = -11746
$ws.Range("H122").Value = 1112090.8
$ws.Range("I122").Value = 1235100.9
$ws.Range("K122").Value = 3705302.7
$ws.Range("M122").Value = -3702852.7
$ws.Range("H127").Value = 11537.186
$ws.Range("I127").Value = 1692.6
$ws.Range("K127").Value = 5077.799999999999
$ws.Range("M127").Value = -117.7999999999993
$ws.Range("H129").Value = 768131.4399999999
$ws.Range("I129").Value = 986197.5600000001
$ws.Range("K129").Value = 2958592.68
$ws.Range("M129").Value = -2953592.68
$ws.Range("H131").Value = 6324.4287
$ws.Range("I131").Value = 3654.7
$ws.Range("K131").Value = 10964.1
$ws.Range("M131").Value = -5924.099999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9604.385
$ws.Range("I32").Value = 7306.34
$ws.Range("K32").Value = 7306.34
$ws.Range("M32").Value = -7019.34
$ws.Range("H132").Value = 3856.087
$ws.Range("I132").Value = 3142.7058
$ws.Range("K132").Value = 9428.117400000001
$ws.Range("M132").Value = -6898.117400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 336.6
$ws.Range("J80").Value = 402.58334
$ws.Range("L80").Value = 402.58334
$ws.Range("N80").Value = -2398.58334
$ws.Range("H83").Value = 336.6
$ws.Range("J83").Value = 402.58334
$ws.Range("L83").Value = 2012.9167
$ws.Range("N83").Value = -11996.9167

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 280.58823
$ws.Range("J7").Value = 332.85715
$ws.Range("L7").Value = 332.85715
$ws.Range("N7").Value = -558.85715
$ws.Range("H31").Value = 3876.8462
$ws.Range("I31").Value = 1833.7778
$ws.Range("J31").Value = 4958.4707
$ws.Range("K31").Value = 1833.7778
$ws.Range("L31").Value = 4958.4707
$ws.Range("M31").Value = -1538.7778
$ws.Range("N31").Value = -5548.4707
$ws.Range("H34").Value = 3876.8462
$ws.Range("I34").Value = 1833.7778
$ws.Range("J34").Value = 4958.4707
$ws.Range("K34").Value = 1833.7778
$ws.Range("L34").Value = 4958.4707
$ws.Range("M34").Value = -1631.7778
$ws.Range("N34").Value = -5362.4707
$ws.Range("H99").Value = 10221.4
$ws.Range("J99").Value = 9690.444
$ws.Range("L99").Value = 9690.444
$ws.Range("N99").Value = -12686.444
$ws.Range("H126").Value = 10221.4
$ws.Range("J126").Value = 9690.444
$ws.Range("L126").Value = 29071.332
$ws.Range("N126").Value = -34011.33199999999
$ws.Range("H132").Value = 2738.5264
$ws.Range("I132").Value = 2403.2666
$ws.Range("K132").Value = 7209.7998
$ws.Range("M132").Value = -4679.7998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1098.1111
$ws.Range("J86").Value = 1088.625
$ws.Range("L86").Value = 3265.875
$ws.Range("N86").Value = -5637.875
$ws.Range("H89").Value = 1098.1111
$ws.Range("J89").Value = 1088.625
$ws.Range("L89").Value = 9797.625
$ws.Range("N89").Value = -21653.625
$ws.Range("H114").Value = 639.75
$ws.Range("I114").Value = 294.1111
$ws.Range("J114").Value = 1084.1428
$ws.Range("K114").Value = 882.3333
$ws.Range("L114").Value = 3252.4284
$ws.Range("M114").Value = 2371.6667
$ws.Range("N114").Value = -9760.428400000001
$ws.Range("H129").Value = 2267.3845
$ws.Range("I129").Value = 1499.5
$ws.Range("J129").Value = 2407
$ws.Range("K129").Value = 4498.5
$ws.Range("L129").Value = 7221
$ws.Range("M129").Value = 501.5
$ws.Range("N129").Value = -17221

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 33447216
$ws.Range("I80").Value = 280254
$ws.Range("J80").Value = 55558524
$ws.Range("K80").Value = 280254
$ws.Range("L80").Value = 55558524
$ws.Range("M80").Value = -279256
$ws.Range("N80").Value = -55560520
$ws.Range("H83").Value = 33447216
$ws.Range("I83").Value = 280254
$ws.Range("J83").Value = 55558524
$ws.Range("K83").Value = 1401270
$ws.Range("L83").Value = 277792620
$ws.Range("M83").Value = -1396278
$ws.Range("N83").Value = -277802604
$ws.Range("H102").Value = 2851.842
$ws.Range("I102").Value = 2005.8
$ws.Range("J102").Value = 6024.5
$ws.Range("K102").Value = 2005.8
$ws.Range("L102").Value = 6024.5
$ws.Range("M102").Value = -383.8
$ws.Range("N102").Value = -9268.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 314816.2
$ws.Range("I93").Value = 2353.5789
$ws.Range("K93").Value = 2353.5789
$ws.Range("M93").Value = -1105.5789
$ws.Range("H132").Value = 4102.241
$ws.Range("I132").Value = 3268.3157
$ws.Range("K132").Value = 9804.947100000001
$ws.Range("M132").Value = -7274.947100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 23817768
$ws.Range("I81").Value = 11106.6
$ws.Range("J81").Value = 83334424
$ws.Range("K81").Value = 22213.2
$ws.Range("L81").Value = 166668848
$ws.Range("M81").Value = -21152.2
$ws.Range("N81").Value = -166670970
$ws.Range("H84").Value = 23817768
$ws.Range("I84").Value = 11106.6
$ws.Range("J84").Value = 83334424
$ws.Range("K84").Value = 111066
$ws.Range("L84").Value = 833344240
$ws.Range("M84").Value = -105762
$ws.Range("N84").Value = -833354848
$ws.Range("H132").Value = 2559.2307
$ws.Range("I132").Value = 2042.9524
$ws.Range("J132").Value = 3161.5557
$ws.Range("K132").Value = 6128.857199999999
$ws.Range("L132").Value = 9484.667099999999
$ws.Range("M132").Value = -3598.857199999999
$ws.Range("N132").Value = -14544.6671
